$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name) to "resources"
$ws.Name = "resources"

# Update resource labels and values
$ws.Range("A2").Value = "PLA_virgin"
$ws.Range("B2").Value = 0.534

$ws.Range("A3").Value = "PLA_recycled"
$ws.Range("B3").Value = 0.342

$ws.Range("A4").Value = "PLA_recycled_industrial"
$ws.Range("B4").Value = 0.471
